# Scheduled runner update: refresh market price / profit data across the
# per-job leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3371.8
$ws.Range("I40").Value = 2700
$ws.Range("J40").Value = 3819.6667
$ws.Range("K40").Value = 2700
$ws.Range("L40").Value = 3819.6667
$ws.Range("M40").Value = -2525
$ws.Range("N40").Value = -4169.6667
$ws.Range("H116").Value = 9458.333000000001
$ws.Range("I116").Value = 12290.5
$ws.Range("K116").Value = 12290.5
$ws.Range("M116").Value = -8848.5
$ws.Range("H137").Value = 4106.3076
$ws.Range("I137").Value = 2772.7827
$ws.Range("J137").Value = 14330
$ws.Range("K137").Value = 8318.348100000001
$ws.Range("L137").Value = 42990
$ws.Range("M137").Value = -5768.348100000001
$ws.Range("N137").Value = -48090
$ws.Range("H138").Value = 2219.4211
$ws.Range("I138").Value = 1830.0385
$ws.Range("J138").Value = 2546
$ws.Range("K138").Value = 5490.1155
$ws.Range("L138").Value = 7638
$ws.Range("M138").Value = -350.1154999999999
$ws.Range("N138").Value = -17918
$ws.Range("H141").Value = 5231.1665
$ws.Range("I141").Value = 2634.0908
$ws.Range("K141").Value = 7902.2724
$ws.Range("M141").Value = -2722.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 51278
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 51278
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 51278
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -51820
$ws.Range("H61").Value = 3541.739
$ws.Range("I61").Value = 3411.5557
$ws.Range("K61").Value = 3411.5557
$ws.Range("M61").Value = -3199.5557
$ws.Range("H74").Value = 2040
$ws.Range("I74").Value = 1296.8889
$ws.Range("J74").Value = 3154.6667
$ws.Range("K74").Value = 1296.8889
$ws.Range("L74").Value = 3154.6667
$ws.Range("M74").Value = -422.8888999999999
$ws.Range("N74").Value = -4902.6667
$ws.Range("H77").Value = 2040
$ws.Range("I77").Value = 1296.8889
$ws.Range("J77").Value = 3154.6667
$ws.Range("K77").Value = 6484.4445
$ws.Range("L77").Value = 15773.3335
$ws.Range("M77").Value = -2116.4445
$ws.Range("N77").Value = -24509.3335
$ws.Range("H132").Value = 5565.5
$ws.Range("I132").Value = 6293.6
$ws.Range("K132").Value = 18880.8
$ws.Range("M132").Value = -16350.8
$ws.Range("H136").Value = 3541.739
$ws.Range("I136").Value = 3411.5557
$ws.Range("K136").Value = 10234.6671
$ws.Range("M136").Value = -7684.667099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2758.9
$ws.Range("I134").Value = 2136.3333
$ws.Range("J134").Value = 3692.75
$ws.Range("K134").Value = 6408.999899999999
$ws.Range("L134").Value = 11078.25
$ws.Range("M134").Value = -3873.999899999999
$ws.Range("N134").Value = -16148.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5255.8477
$ws.Range("I31").Value = 1130.5667
$ws.Range("J31").Value = 9523.379000000001
$ws.Range("K31").Value = 1130.5667
$ws.Range("L31").Value = 9523.379000000001
$ws.Range("M31").Value = -835.5667000000001
$ws.Range("N31").Value = -10113.379
$ws.Range("H34").Value = 5255.8477
$ws.Range("I34").Value = 1130.5667
$ws.Range("J34").Value = 9523.379000000001
$ws.Range("K34").Value = 1130.5667
$ws.Range("L34").Value = 9523.379000000001
$ws.Range("M34").Value = -928.5667000000001
$ws.Range("N34").Value = -9927.379000000001
$ws.Range("H58").Value = 1846.8889
$ws.Range("I58").Value = 1592
$ws.Range("J58").Value = 2050.8
$ws.Range("K58").Value = 1592
$ws.Range("L58").Value = 2050.8
$ws.Range("M58").Value = -1389
$ws.Range("N58").Value = -2456.8
$ws.Range("H99").Value = 1761.28
$ws.Range("I99").Value = 1166.6666
$ws.Range("J99").Value = 1842.3636
$ws.Range("K99").Value = 1166.6666
$ws.Range("L99").Value = 1842.3636
$ws.Range("M99").Value = 331.3334
$ws.Range("N99").Value = -4838.3636
$ws.Range("H107").Value = 3473118
$ws.Range("I107").Value = 4808498
$ws.Range("J107").Value = 1130
$ws.Range("K107").Value = 4808498
$ws.Range("L107").Value = 1130
$ws.Range("M107").Value = -4806578
$ws.Range("N107").Value = -4970
$ws.Range("H126").Value = 1761.28
$ws.Range("I126").Value = 1166.6666
$ws.Range("J126").Value = 1842.3636
$ws.Range("K126").Value = 3499.9998
$ws.Range("L126").Value = 5527.0908
$ws.Range("M126").Value = -1029.9998
$ws.Range("N126").Value = -10467.0908
$ws.Range("H132").Value = 8335562
$ws.Range("I132").Value = 1879.4375
$ws.Range("J132").Value = 41670292
$ws.Range("K132").Value = 5638.3125
$ws.Range("L132").Value = 125010876
$ws.Range("M132").Value = -3108.3125
$ws.Range("N132").Value = -125015936
$ws.Range("H136").Value = 1846.8889
$ws.Range("I136").Value = 1592
$ws.Range("J136").Value = 2050.8
$ws.Range("K136").Value = 4776
$ws.Range("L136").Value = 6152.400000000001
$ws.Range("M136").Value = -2226
$ws.Range("N136").Value = -11252.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 749.5
$ws.Range("I5").Value = 752.8
$ws.Range("K5").Value = 2258.4
$ws.Range("M5").Value = -2146.4
$ws.Range("H135").Value = 749.5
$ws.Range("I135").Value = 752.8
$ws.Range("K135").Value = 6775.2
$ws.Range("M135").Value = -4240.2
$ws.Range("H140").Value = 1934.7059
$ws.Range("I140").Value = 1681.8182
$ws.Range("J140").Value = 2398.3333
$ws.Range("K140").Value = 5045.4546
$ws.Range("L140").Value = 7194.999899999999
$ws.Range("M140").Value = 134.5454
$ws.Range("N140").Value = -17554.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 90009
$ws.Range("J25").Value = 90009
$ws.Range("L25").Value = 90009
$ws.Range("N25").Value = -91067
$ws.Range("H126").Value = 2009.2142
$ws.Range("I126").Value = 2038.2
$ws.Range("J126").Value = 1936.75
$ws.Range("K126").Value = 6114.6
$ws.Range("L126").Value = 5810.25
$ws.Range("M126").Value = -3644.6
$ws.Range("N126").Value = -10750.25
$ws.Range("H132").Value = 2846.6538
$ws.Range("I132").Value = 2511.6875
$ws.Range("J132").Value = 3382.6
$ws.Range("K132").Value = 7535.0625
$ws.Range("L132").Value = 10147.8
$ws.Range("M132").Value = -5005.0625
$ws.Range("N132").Value = -15207.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 64666.668
$ws.Range("J106").Value = 64666.668
$ws.Range("L106").Value = 64666.668
$ws.Range("N106").Value = -67190.66800000001
$ws.Range("H132").Value = 2019.3914
$ws.Range("I132").Value = 1579.2424
$ws.Range("J132").Value = 3136.6924
$ws.Range("K132").Value = 4737.7272
$ws.Range("L132").Value = 9410.0772
$ws.Range("M132").Value = -2207.7272
$ws.Range("N132").Value = -14470.0772
$ws.Range("H136").Value = 6946459
$ws.Range("J136").Value = 12822402
$ws.Range("L136").Value = 38467206
$ws.Range("N136").Value = -38472306

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 70750
$ws.Range("J105").Value = 70750
$ws.Range("L105").Value = 70750
$ws.Range("N105").Value = -77738
$ws.Range("H107").Value = 760.6667
$ws.Range("I107").Value = 302
$ws.Range("J107").Value = 990
$ws.Range("K107").Value = 906
$ws.Range("L107").Value = 2970
$ws.Range("M107").Value = 1014
$ws.Range("N107").Value = -6810
$ws.Range("H132").Value = 7939548
$ws.Range("I132").Value = 3158.9285
$ws.Range("J132").Value = 23812326
$ws.Range("K132").Value = 9476.7855
$ws.Range("L132").Value = 71436978
$ws.Range("M132").Value = -6946.7855
$ws.Range("N132").Value = -71442038
$ws.Range("H136").Value = 2897.1052
$ws.Range("I136").Value = 2895.8696
$ws.Range("K136").Value = 8687.6088
$ws.Range("M136").Value = -6137.6088
